$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1: B1 = "Хэш" (new), C1 = "Время обработки" (new), D1 = "Хэммингово расстояние" (moved from old B1)
$ws.Range("B1").Value = "Хэш"
$ws.Range("C1").Value = "Время обработки"
$ws.Range("D1").Value = "Хэммингово расстояние"

# Force column B (rows 2-63) to text so 68-character 0/1 hash strings are preserved exactly (not converted to scientific-notation numbers)
$ws.Range("B2:B63").NumberFormat = "@"

$rows = @(
    @{ Row = 2; B = "1010010000000000100000100000000010000110000000001000000000000000"; C = 0.013809; D = 10 },
    @{ Row = 3; B = "1110010000000000100000100000000010000110000000001000000000000000"; C = 0; D = 9 },
    @{ Row = 4; B = "1010000000000000100000000000000000000000000000000000000000000000"; C = 0.013875; D = 10 },
    @{ Row = 5; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0.00719; D = 0 },
    @{ Row = 6; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0; D = 0 },
    @{ Row = 7; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0; D = 0 },
    @{ Row = 8; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0; D = 0 },
    @{ Row = 9; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0; D = 0 },
    @{ Row = 10; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0.006722; D = 0 },
    @{ Row = 11; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0; D = 0 },
    @{ Row = 12; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0; D = 0 },
    @{ Row = 13; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0; D = 0 },
    @{ Row = 14; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0; D = 0 },
    @{ Row = 15; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0.006877; D = 0 },
    @{ Row = 16; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0; D = 0 },
    @{ Row = 17; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0; D = 0 },
    @{ Row = 18; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0; D = 0 },
    @{ Row = 19; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0; D = 0 },
    @{ Row = 20; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0.007045; D = 0 },
    @{ Row = 21; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0; D = 0 },
    @{ Row = 22; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0; D = 0 },
    @{ Row = 23; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0.007124; D = 0 },
    @{ Row = 24; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0.007009; D = 0 },
    @{ Row = 25; B = "1110101100000000101000001000000010000100000000001000001000000000"; C = 0; D = 0 },
    @{ Row = 26; B = "1110101100000000101000001000000010000000000000001000001000000000"; C = 0.000746; D = 1 },
    @{ Row = 27; B = "1110101100000000101000001000000010000000000000001000001000000000"; C = 0.007046; D = 1 },
    @{ Row = 28; B = "1110101100000000101000001000000010000000000000001000001000000000"; C = 0.006787; D = 1 },
    @{ Row = 29; B = "1110101100000000101000001000000010000000000000001000001000000000"; C = 0.006239; D = 1 },
    @{ Row = 30; B = "1110101100000000101000001000000010000000000000001000001000000000"; C = 0.000499; D = 1 },
    @{ Row = 31; B = "1110101100000000101000001000000010000000000000001000001000000000"; C = 0.006919; D = 1 },
    @{ Row = 32; B = "1110101100000000101000001000000010000000000000001000001000000000"; C = 0.006982; D = 1 },
    @{ Row = 33; B = "1110101100000000101000001000000010000000000000001000000000000000"; C = 0.00706; D = 2 },
    @{ Row = 34; B = "1010000000000000100000000000000000000000000000000000000000000000"; C = 0; D = 10 },
    @{ Row = 35; B = "1110101101000000101001000001000010100100000000001000011000000000"; C = 0; D = 6 },
    @{ Row = 36; B = "1110100100000000100000000000000010000010100000000000111000000000"; C = 0.006895; D = 9 },
    @{ Row = 37; B = "1110100000000000100100001000000000000000000000000000000001000000"; C = 0; D = 9 },
    @{ Row = 38; B = "1110110001100000101000001000100110000000000000001000000000000000"; C = 0.00694; D = 9 },
    @{ Row = 39; B = "1011100000000000110000000000000000000000000000000000000000000000"; C = 0; D = 11 },
    @{ Row = 40; B = "1010100100000000100000000000000010000010000000001000000000000000"; C = 0.006955; D = 7 },
    @{ Row = 41; B = "1011100000000000110000001000100000000000100000000000000000000000"; C = 0; D = 12 },
    @{ Row = 42; B = "1011000000000000100000000000000000000000000000000000000000000000"; C = 0.007021; D = 11 },
    @{ Row = 43; B = "1011100100000000110000000000000010000000000000000001000000000000"; C = 0.006922; D = 10 },
    @{ Row = 44; B = "1011000000110000110000001100000010000000000000000000000000000000"; C = 0; D = 13 },
    @{ Row = 45; B = "1011101100000000110000000000000010000000000000000001000000000000"; C = 0.006956; D = 9 },
    @{ Row = 46; B = "1011100100000000110001001000000000010000000000000000000000000000"; C = 0.009098; D = 11 },
    @{ Row = 47; B = "1011100000000000110000001000000000000000000000000000000000000000"; C = 0.00472; D = 10 },
    @{ Row = 48; B = "1011001100000000110000000000000010000000000000000000000000000000"; C = 0.006929; D = 9 },
    @{ Row = 49; B = "1111100100000000100001000001000010000010000000001000000000000000"; C = 0.006958; D = 9 },
    @{ Row = 50; B = "1110100100000000110011000000000000000000000000000010000000000000"; C = 0.006953; D = 11 },
    @{ Row = 51; B = "1011101100000000110000001000000000000000000000000001000000000000"; C = 0.006955; D = 9 },
    @{ Row = 52; B = "1011101100000000110000001000000000000000000000000001000000000000"; C = 0.006928; D = 9 },
    @{ Row = 53; B = "1011100100000000110000000000000000000000000000000001000000000000"; C = 0.007033; D = 11 },
    @{ Row = 54; B = "1010001000000000110000000000000000000000000000000000000000000000"; C = 0.013799; D = 10 },
    @{ Row = 55; B = "1010000000000000100000000000000000010000000000000000000000000000"; C = 0.006928; D = 11 },
    @{ Row = 56; B = "1100111000010000101100000000000000000000100000000000000000000000"; C = 0.013996; D = 11 },
    @{ Row = 57; B = "1011000000000000100000100000000000000000000000001000000000000000"; C = 0.006952; D = 11 },
    @{ Row = 58; B = "1011101100000000110001000000000000010010000000000001000000000000"; C = 0.01379; D = 13 },
    @{ Row = 59; B = "1010000000000000100000100000000010001000000000001000000000000000"; C = 0.028726; D = 10 },
    @{ Row = 60; B = "1010100000010000110000001100000000000010100000011000000000010000"; C = 0.033469; D = 14 },
    @{ Row = 61; B = "1011101001000000111000000000000000000000100000000000010100000000"; C = 0.041659; D = 13 },
    @{ Row = 62; B = "1011101010000000110000000000000011000000000000000000001000010000"; C = 0.041534; D = 11 },
    @{ Row = 63; B = "1011000000000000100000000000000000000000000000000000000000000000"; C = 0.048903; D = 11 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
